# Worked on Time Series for RCH and EVT in MODFLOW 6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Actual (hours)" values (column I) reflecting the new work done.
$ws.Range("I3").Value = 50     # Time Lists
$ws.Range("I33").Value = 2     # RCH
$ws.Range("I34").Value = 1     # EVT
$ws.Range("I37").Value = 2     # Import model results

# Update the active selection to match where the user ended up working (K35).
$ws.Range("K35").Select()
